$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3169.3333
$ws.Range("I74").Value = 3100
$ws.Range("J74").Value = 3183.2
$ws.Range("K74").Value = 3100
$ws.Range("L74").Value = 3183.2
$ws.Range("M74").Value = -2164
$ws.Range("N74").Value = -5055.2
$ws.Range("H77").Value = 3169.3333
$ws.Range("I77").Value = 3100
$ws.Range("J77").Value = 3183.2
$ws.Range("K77").Value = 15500
$ws.Range("L77").Value = 15916
$ws.Range("M77").Value = -10820
$ws.Range("N77").Value = -25276
$ws.Range("H86").Value = 157191710
$ws.Range("I86").Value = 157191710
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 157191710
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -157190587
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 157191710
$ws.Range("I89").Value = 157191710
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 785958550
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -785952934
$ws.Range("N89").ClearContents()
$ws.Range("H137").Value = 1313.6774
$ws.Range("I137").Value = 954.3
$ws.Range("K137").Value = 2862.9
$ws.Range("M137").Value = -312.8999999999996
$ws.Range("H141").Value = 3752.6365
$ws.Range("I141").Value = 1915.1666
$ws.Range("J141").Value = 12021.25
$ws.Range("K141").Value = 5745.4998
$ws.Range("L141").Value = 36063.75
$ws.Range("M141").Value = -565.4997999999996
$ws.Range("N141").Value = -46423.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2243.7222
$ws.Range("I45").Value = 2150
$ws.Range("J45").Value = 2360.875
$ws.Range("K45").Value = 2150
$ws.Range("L45").Value = 2360.875
$ws.Range("M45").Value = -1773
$ws.Range("N45").Value = -3114.875
$ws.Range("H46").Value = 16669333
$ws.Range("J46").Value = 16669333
$ws.Range("L46").Value = 16669333
$ws.Range("N46").Value = -16669971
$ws.Range("H61").Value = 4184.1763
$ws.Range("I61").Value = 3433.5386
$ws.Range("J61").Value = 4648.857
$ws.Range("K61").Value = 3433.5386
$ws.Range("L61").Value = 4648.857
$ws.Range("M61").Value = -3221.5386
$ws.Range("N61").Value = -5072.857
$ws.Range("H74").Value = 2443.7708
$ws.Range("I74").Value = 1416.3793
$ws.Range("J74").Value = 4011.8948
$ws.Range("K74").Value = 1416.3793
$ws.Range("L74").Value = 4011.8948
$ws.Range("M74").Value = -542.3793000000001
$ws.Range("N74").Value = -5759.8948
$ws.Range("H77").Value = 2443.7708
$ws.Range("I77").Value = 1416.3793
$ws.Range("J77").Value = 4011.8948
$ws.Range("K77").Value = 7081.896500000001
$ws.Range("L77").Value = 20059.474
$ws.Range("M77").Value = -2713.896500000001
$ws.Range("N77").Value = -28795.474
$ws.Range("H102").Value = 3250
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3250
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -6494
$ws.Range("H132").Value = 4306.7144
$ws.Range("I132").Value = 2821.7222
$ws.Range("J132").Value = 6979.7
$ws.Range("K132").Value = 8465.1666
$ws.Range("L132").Value = 20939.1
$ws.Range("M132").Value = -5935.1666
$ws.Range("N132").Value = -25999.1
$ws.Range("H136").Value = 4184.1763
$ws.Range("I136").Value = 3433.5386
$ws.Range("J136").Value = 4648.857
$ws.Range("K136").Value = 10300.6158
$ws.Range("L136").Value = 13946.571
$ws.Range("M136").Value = -7750.6158
$ws.Range("N136").Value = -19046.571

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2021
$ws.Range("I99").Value = 1729
$ws.Range("K99").Value = 1729
$ws.Range("M99").Value = -231
$ws.Range("H134").Value = 2146.2444
$ws.Range("I134").Value = 1804.421
$ws.Range("K134").Value = 5413.263
$ws.Range("M134").Value = -2878.263

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1667.4375
$ws.Range("I58").Value = 1138.7391
$ws.Range("K58").Value = 1138.7391
$ws.Range("M58").Value = -935.7391
$ws.Range("H132").Value = 4763814.5
$ws.Range("I132").Value = 1642.2727
$ws.Range("J132").Value = 12822875
$ws.Range("K132").Value = 4926.8181
$ws.Range("L132").Value = 38468625
$ws.Range("M132").Value = -2396.8181
$ws.Range("N132").Value = -38473685
$ws.Range("H134").Value = 6616.577
$ws.Range("I134").Value = 6667.952
$ws.Range("K134").Value = 20003.856
$ws.Range("M134").Value = -17468.856
$ws.Range("H136").Value = 1667.4375
$ws.Range("I136").Value = 1138.7391
$ws.Range("K136").Value = 3416.2173
$ws.Range("M136").Value = -866.2173000000003

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 975.2174
$ws.Range("I20").Value = 810
$ws.Range("K20").Value = 2430
$ws.Range("M20").Value = -2203

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7308.125
$ws.Range("I70").Value = 7816
$ws.Range("J70").Value = 5599.8184
$ws.Range("K70").Value = 7816
$ws.Range("L70").Value = 5599.8184
$ws.Range("M70").Value = -7546
$ws.Range("N70").Value = -6139.8184
$ws.Range("H73").Value = 7308.125
$ws.Range("I73").Value = 7816
$ws.Range("J73").Value = 5599.8184
$ws.Range("K73").Value = 7816
$ws.Range("L73").Value = 5599.8184
$ws.Range("M73").Value = -6880
$ws.Range("N73").Value = -7471.8184
$ws.Range("H123").Value = 8321
$ws.Range("J123").Value = 8321
$ws.Range("L123").Value = 8321
$ws.Range("N123").Value = -13221

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 881.45
$ws.Range("I16").Value = 601.875
$ws.Range("J16").Value = 1999.75
$ws.Range("K16").Value = 601.875
$ws.Range("L16").Value = 1999.75
$ws.Range("M16").Value = -431.875
$ws.Range("N16").Value = -2339.75
$ws.Range("H22").Value = 10785.4
$ws.Range("I22").Value = 890.4
$ws.Range("J22").Value = 20680.4
$ws.Range("K22").Value = 890.4
$ws.Range("L22").Value = 20680.4
$ws.Range("M22").Value = -595.4
$ws.Range("N22").Value = -21270.4
$ws.Range("H27").Value = 10785.4
$ws.Range("I27").Value = 890.4
$ws.Range("J27").Value = 20680.4
$ws.Range("K27").Value = 890.4
$ws.Range("L27").Value = 20680.4
$ws.Range("M27").Value = -783.4
$ws.Range("N27").Value = -20894.4
$ws.Range("H46").Value = 674.6667
$ws.Range("I46").Value = 1179
$ws.Range("J46").Value = 422.5
$ws.Range("K46").Value = 1179
$ws.Range("L46").Value = 422.5
$ws.Range("M46").Value = -991
$ws.Range("N46").Value = -798.5
$ws.Range("H68").Value = 1943.5349
$ws.Range("I68").Value = 1829.7333
$ws.Range("J68").Value = 2206.1538
$ws.Range("K68").Value = 1829.7333
$ws.Range("L68").Value = 2206.1538
$ws.Range("M68").Value = -1080.7333
$ws.Range("N68").Value = -3704.1538
$ws.Range("H71").Value = 1943.5349
$ws.Range("I71").Value = 1829.7333
$ws.Range("J71").Value = 2206.1538
$ws.Range("K71").Value = 9148.666500000001
$ws.Range("L71").Value = 11030.769
$ws.Range("M71").Value = -5404.666500000001
$ws.Range("N71").Value = -18518.769
$ws.Range("H132").Value = 2897.6875
$ws.Range("I132").Value = 2323.1667
$ws.Range("K132").Value = 6969.500100000001
$ws.Range("M132").Value = -4439.500100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3404.6667
$ws.Range("I136").Value = 3328.9565
$ws.Range("K136").Value = 9986.869499999999
$ws.Range("M136").Value = -7436.869499999999
